$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume change (E) columns per latest scrape.
# Cells whose new value looks numeric are pinned to Text format first so Excel
# keeps them as strings (preserving exact formatting such as trailing zeros),
# matching how this sheet stores all Price/Volume figures as text.

$ws.Range("D2").Value = '41.915.04'
$ws.Range("D3").Value = '2.268.90'
$ws.Range("E3").Value = '  +3.04%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.17'
$ws.Range("E5").Value = '  +3.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.06'
$ws.Range("E6").Value = '  +7.92%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.525'
$ws.Range("E7").Value = '  +3.40%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  +4.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.75'
$ws.Range("E10").Value = '  +8.87%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.50'
$ws.Range("E11").Value = '  +9.09%  '
$ws.Range("E12").Value = '  +2.88%  '
$ws.Range("E13").Value = '  +3.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.70'
$ws.Range("E14").Value = '  +3.93%  '
$ws.Range("D15").Value = '2.620.25'
$ws.Range("E15").Value = '  +3.07%  '
$ws.Range("E16").Value = '  +3.71%  '
$ws.Range("D17").Value = '2.268.32'
$ws.Range("E17").Value = '  +2.61%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.756'
$ws.Range("E18").Value = '  +3.88%  '
$ws.Range("D19").Value = '41.807.14'
$ws.Range("E19").Value = '  +5.20%  '
$ws.Range("E20").Value = '  +10.69%  '
$ws.Range("D21").Value = '0.0₃0908'
$ws.Range("E21").Value = '  +3.05%  '
$ws.Range("E22").Value = '  +3.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.36'
$ws.Range("E23").Value = '  +3.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '241.51'
$ws.Range("E24").Value = '  +1.78%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.59'
$ws.Range("E25").Value = '  +6.48%  '
$ws.Range("E27").Value = '  +5.47%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.89'
$ws.Range("E28").Value = '  +1.81%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.25'
$ws.Range("E29").Value = '  +9.42%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.71'
$ws.Range("E30").Value = '  +5.88%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.22'
$ws.Range("E31").Value = '  +9.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '158.18'
$ws.Range("E32").Value = '  +1.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.19'
$ws.Range("E34").Value = '  +5.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0741'
$ws.Range("E35").Value = '  +5.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.07'
$ws.Range("E36").Value = '  +6.88%  '
$ws.Range("E37").Value = '  +3.26%  '
$ws.Range("E38").Value = '  +7.90%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '16.55'
$ws.Range("E39").Value = '  +9.11%  '
$ws.Range("E40").Value = '  +3.41%  '
$ws.Range("E41").Value = '  +6.82%  '
$ws.Range("E42").Value = '  +6.86%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.42'
$ws.Range("E43").Value = '  +16.99%  '
$ws.Range("D44").Value = '2.049.61'
$ws.Range("E44").Value = '  -3.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0279'
$ws.Range("E46").Value = '  +3.39%  '
$ws.Range("E47").Value = '  +10.27%  '
$ws.Range("E48").Value = '  -3.97%  '
$ws.Range("D49").Value = '2.491.94'
$ws.Range("E49").Value = '  +3.29%  '
$ws.Range("E50").Value = '  +3.24%  '
$ws.Range("E51").Value = '  +4.40%  '
